$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: adjust formulas
$ws.Range("B15").Formula = "=SUM(B5:B14)"
$ws.Range("C15").Formula = "=SUM(C5:C13)"
$ws.Range("D15").Formula = "=SUM(D5:D13)"
$ws.Range("E15:H15").FormulaR1C1 = "=SUM(R5C:R13C)"

# New row 17: MINIMO
$ws.Range("A17").Value = "MÍNIMO"
$ws.Range("B17:G17").FormulaR1C1 = "=MIN(R5C:R13C)"

# New row 18: MAXIMO
$ws.Range("A18").Value = "MÁXIMO"
$ws.Range("B18:G18").FormulaR1C1 = "=MAX(R5C:R13C)"

# New row 19: MAXIMO_6MESES
$ws.Range("A19").Value = "MÁXIMO_6MESES"
$ws.Range("B19").Formula = "=MAX(B5:G13)"

# New row 20: PROMEDIO
$ws.Range("A20").Value = "PROMEDIO"
$ws.Range("B20:G20").FormulaR1C1 = "=AVERAGE(R5C:R13C)"

# New row 21: CONTAR
$ws.Range("A21").Value = "CONTAR"
$ws.Range("B21:G21").FormulaR1C1 = "=COUNT(R5C:R13C)"

# New row 22: CONTAR_GENERAL
$ws.Range("A22").Value = "CONTAR_GENERAL"
$ws.Range("B22").Formula = "=COUNT(B5:G13)"

# Apply currency style (s=2) to new numeric rows 17,18,20 matching existing row 15 styling
$styleSource = $ws.Range("B15")
$ws.Range("B17:G17").NumberFormat = $styleSource.NumberFormat
$ws.Range("B18:G18").NumberFormat = $styleSource.NumberFormat
$ws.Range("B19").NumberFormat = $styleSource.NumberFormat
$ws.Range("B20:G20").NumberFormat = $styleSource.NumberFormat

$ws.Calculate()

$ws.Range("A2").Select()
$ws.Range("B22").Select()
